# UGN-356 - add "date" column (prop dateFormat support)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "date" column
$ws.Range("C1").Value = "date"
$ws.Range("C1").Font.Color = 0

# First data row gets the value + the explicit date number format + font
$ws.Range("C2").Value = 43893
$ws.Range("C2").Font.Color = 0
$ws.Range("C2").NumberFormat = "mm-dd-yy"

# Remaining rows: copy the same format down (so the style is reused,
# matching the single shared cellXfs entry), then (re)apply their values
$ws.Range("C3").Value = 40272
$ws.Range("C4").Value = 34392

$ws.Range("C2").Copy()
$ws.Range("C3:C4").PasteSpecial(-4122)

$ws.Range("C3").Value = 40272
$ws.Range("C4").Value = 34392

$null = $ws.Range("C4").Select()
